$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the Comment for designator C7 from "5nF" to "20nF".
#    A leading apostrophe keeps this text cell's existing "stored as text"
#    (quote-prefix) formatting, same as the other numeric-looking Comment
#    values in this column (3k, 100k, 1k, 100nF, 1uF, ...).
$ws.Range("B8").Value = "'20nF"

# 2) Rows 12-16 (U1..U5/U8 footprint rows) previously used a cell style
#    that explicitly applied a (no-op/"none") fill on top of the font+border
#    formatting. Re-apply "no fill" so these cells match the plain
#    font+border style already used by rows 9-11, removing the redundant
#    fill-applying style.
$ws.Range("A12:B16").Interior.Pattern = -4142

# 3) Update the saved cursor/selection position to B9
$ws.Range("B9").Select()
